$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44650
$ws.Range("K2").Value = "Angeleno"
$ws.Range("N2").Value = 17000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 17500
$ws.Range("S2").Value = 972

# Row 3
$ws.Range("D3").Value = 44285
$ws.Range("K3").Value = "Angeleno"
$ws.Range("L3").Value = "Primera"

# Row 4
$ws.Range("D4").Value = 44278
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 15500
$ws.Range("S4").Value = 861

# Row 5
$ws.Range("D5").Value = 44614
$ws.Range("K5").Value = "Angeleno"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 19000
$ws.Range("P5").Value = 18500
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 1028

# Row 6
$ws.Range("D6").Value = 44169
$ws.Range("L6").Value = "Tercera"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("S6").Value = 1361

# Row 7
$ws.Range("D7").Value = 44217
$ws.Range("K7").Value = "Black Amber"
$ws.Range("N7").Value = 16000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 16500
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 917

# Row 8
$ws.Range("D8").Value = 44596
$ws.Range("M8").Value = 250
$ws.Range("Q8").Value = "`$/caja 18 kilos granel"

# Row 9
$ws.Range("D9").Value = 44706
$ws.Range("K9").Value = "Angeleno"
$ws.Range("M9").Value = 300
$ws.Range("Q9").Value = "`$/bandeja 18 kilos granel"

# Row 10
$ws.Range("D10").Value = 44580
$ws.Range("K10").Value = "Black Amber"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 270
$ws.Range("N10").Value = 19000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 19500
$ws.Range("R10").Value = "Región Metropolitana"
$ws.Range("S10").Value = 1083

# Row 11
$ws.Range("D11").Value = 44238
$ws.Range("L11").Value = "Segunda"
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 14500
$ws.Range("S11").Value = 806

# Row 12
$ws.Range("D12").Value = 44238
$ws.Range("K12").Value = "Fortuna"
$ws.Range("M12").Value = 300

# Row 13
$ws.Range("D13").Value = 44174
$ws.Range("K13").Value = "Angeleno"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 270
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 20500
$ws.Range("S13").Value = 1139

# Row 14
$ws.Range("D14").Value = 44921
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 450
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 19111
$ws.Range("Q14").Value = "`$/caja 18 kilos granel"
$ws.Range("S14").Value = 1062

# Row 16
$ws.Range("D16").Value = 44314
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 14500
$ws.Range("Q16").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S16").Value = 806

# Row 17
$ws.Range("D17").Value = 44245
$ws.Range("K17").Value = "Black Amber"
$ws.Range("L17").Value = "Primera"
$ws.Range("N17").Value = 14000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 14500
$ws.Range("R17").Value = "Región de O'Higgins"
$ws.Range("S17").Value = 806

# Row 18
$ws.Range("D18").Value = 44229
$ws.Range("K18").Value = "Fortuna"
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 14000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 14500
$ws.Range("Q18").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S18").Value = 806

# Row 19
$ws.Range("D19").Value = 44587
$ws.Range("K19").Value = "Black Amber"
$ws.Range("Q19").Value = "`$/caja 18 kilos granel"

# Row 20
$ws.Range("D20").Value = 44239
$ws.Range("K20").Value = "Fortuna"
$ws.Range("L20").Value = "Primera"
$ws.Range("N20").Value = 15000
$ws.Range("O20").Value = 16000
$ws.Range("P20").Value = 15500
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 861

# Row 21
$ws.Range("D21").Value = 44175
$ws.Range("K21").Value = "Angeleno"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 21000
$ws.Range("O21").Value = 22000
$ws.Range("P21").Value = 21500
$ws.Range("S21").Value = 1194

# Row 22
$ws.Range("D22").Value = 44574
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 18000
$ws.Range("O22").Value = 19000
$ws.Range("P22").Value = 18500
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 1028

# Row 23
$ws.Range("D23").Value = 44628
$ws.Range("K23").Value = "Black Amber"
$ws.Range("M23").Value = 270
$ws.Range("N23").Value = 15000
$ws.Range("O23").Value = 16000
$ws.Range("P23").Value = 15500
$ws.Range("S23").Value = 861
